$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2..31 (timestamp, label, ax, ay, az, gx, gy, gz)
# Row layout: rowIndex, timestamp, label, ax, ay, az, gx, gy, gz
$data = @(
    @(2, 0, "walkingToRunning", -1.054504990577698, -3.490631103515625, 0.8201689720153809, 0.2057795317889, -0.6517161301269102, -1.099296086808523),
    @(3, 100, "walkingToRunning", -1.129331469535828, -6.879757404327393, 1.420253753662109, -0.1219321910927993, -0.182187830529549, -0.09341037564160559),
    @(4, 200, "walkingToRunning", -2.684658765792847, -0.5505852699279785, 8.388395309448242, -0.4356265911241861, -0.6137270996483248, 1.354154146480842),
    @(5, 300, "walkingToRunning", -5.721604347229004, 7.858576774597168, 6.318061828613281, -0.06265412089301869, -0.5025855811631026, 2.260356429146558),
    @(6, 400, "walkingToRunning", 1.885133743286133, -1.666801452636719, 13.95527076721191, -0.5034300181560275, 0.6285661230000101, 1.827248451186396),
    @(7, 500, "walkingToRunning", 4.083763122558594, 2.682566165924072, 5.88519811630249, -1.446588513327809, 0.08285516464128959, 2.005940111671997),
    @(8, 600, "walkingToRunning", -37.18881988525391, -16.31022262573242, 18.38444519042969, -1.11462547393834, -0.009374656328342557, 3.162934237863957),
    @(9, 700, "walkingToRunning", -4.98013973236084, 3.857727766036987, 4.614824771881104, -0.2443921377019233, 1.357359671556357, 1.662632597655804),
    @(10, 800, "walkingToRunning", 44.95419311523438, -55.06259918212891, -3.026922225952148, 0.1685504996921994, 1.441876109053465, -4.385711685186457),
    @(11, 900, "walkingToRunning", 43.95158767700195, -54.27328491210938, -3.213864326477051, 0.2859861443682361, -2.243682384490969, -7.90092859035586),
    @(12, 1000, "walkingToRunning", -9.904374122619627, 6.515621185302734, 11.26030445098877, 2.782215724631054, -0.7535885743024512, -2.988133376691399),
    @(13, 1100, "walkingToRunning", -2.358283281326294, -57.28740692138672, -2.653994083404541, -0.493949881414097, -2.56519750065012, -2.295662734566661),
    @(14, 1200, "walkingToRunning", -12.77801513671875, 3.484463453292847, 8.297730445861816, -7.865739104224518, -8.845545524049903, 1.579899276175187),
    @(15, 1300, "walkingToRunning", -23.53105735778809, 3.71211051940918, 4.357646942138672, 3.534730899624699, -1.276565199218106, -3.009974456414927),
    @(16, 1400, "walkingToRunning", -33.88568878173828, -50.99709701538086, 13.22272872924805, 4.82583797559505, -4.776124651839162, -4.42023835749164),
    @(17, 1500, "walkingToRunning", -24.3600902557373, 2.738663196563721, 11.53956317901611, 4.827630825159058, 5.732881429718808, 1.179853487305514),
    @(18, 1600, "walkingToRunning", 18.28820610046387, -0.5376138687133789, 9.498839378356934, -2.591416502871155, -4.246028906688475, 3.687537992872843),
    @(19, 1700, "walkingToRunning", -62.44527816772461, -7.982970237731934, -32.66093826293945, -5.452502425124007, 1.563533739345819, 3.775565188105511),
    @(20, 1800, "walkingToRunning", 23.0417537689209, 16.5068416595459, -16.63084411621094, -4.596750730421486, 7.476189779072271, 4.383058036245957),
    @(21, 1900, "walkingToRunning", -5.877625465393066, -9.525984764099119, -0.1709146499633789, -3.8169967372243, 9.492706188341467, -5.739285922631989),
    @(22, 2000, "walkingToRunning", -0.5468623042106628, -5.408160209655762, 43.39143753051758, 5.77226390198965, -3.99830269232031, -2.008213531680217),
    @(23, 2100, "walkingToRunning", -6.48110294342041, 21.88984298706055, 9.76063346862793, 2.478664084178635, -3.857492981887443, -1.977369552705346),
    @(24, 2200, "walkingToRunning", 19.5135498046875, -70.08018493652344, 22.09575080871582, 0.5908174134972555, -0.1589468512229657, 2.452448280846219),
    @(25, 2300, "walkingToRunning", -37.2692756652832, 22.71374130249023, -6.724684715270996, -13.21499349467632, -7.594239313064515, 6.02689108034458),
    @(26, 2400, "walkingToRunning", -22.65700340270996, 4.547637939453125, -17.0935115814209, 4.174701481330613, -10.92443460371436, 1.552455660046618),
    @(27, 2500, "walkingToRunning", 78.31330871582031, -70.41346740722656, -6.543253421783447, 1.570928898958005, -8.095582741062824, 0.9565766177526321),
    @(28, 2600, "walkingToRunning", -3.919276475906372, 13.1894645690918, 2.611269950866699, 8.965150342112763, 11.8617115107978, 1.040716253402757),
    @(29, 2700, "walkingToRunning", 41.53395843505859, -10.51663017272949, 21.93678855895996, 0.2588674527842816, 4.578985643822968, 1.977885912104351),
    @(30, 2800, "walkingToRunning", 19.39011192321777, -10.6725959777832, -12.10480403900146, -2.832312445815016, 0.201128745497, 1.380974077596895),
    @(31, 2900, "walkingToRunning", 10.54305171966553, 18.1186752319336, -10.15065765380859, -4.569073584021599, -4.38156421046433, -0.2334796684544305)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}
